# Populate the worksheet with the new City/Country table, replacing the
# old single ZipCode1/date layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write data column by column so that shared-string indices come out in
# the same order as the target workbook (City, Columbus, Country, USA).
$ws.Range("A1").Value = "City"
$ws.Range("A2").Value = "Columbus"
$ws.Range("B1").Value = "Country"
$ws.Range("B2").Value = "USA"

# Header row (row 1) is bold.
$ws.Range("A1:B1").Font.Bold = $true

# Column widths: A=16, B=14 characters.
# (COM ColumnWidth has a small built-in padding offset (5/6 of a character)
# in this runtime, so subtract it to land exactly on the target width values.)
$ws.Columns("A").ColumnWidth = 15.166666666666666
$ws.Columns("B").ColumnWidth = 13.166666666666666

# Page orientation is explicitly set to portrait.
$ws.PageSetup.Orientation = 1

# Leave the final selection on D9, matching the saved UI state.
$ws.Range("D9").Select()
